# Applies the "#5: insurance, claim, debt, investment done" edit.
# Adds extra metadata/lookup columns (property_category, category, date,
# legislator_name, legislator_id, source_file, index, ...) to the
# "insurance" (sheet7), "claim" (sheet8) and "debt" (sheet9) worksheets,
# mirroring the layout already used on the other asset-category sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 7: 保險 (insurance)
# ---------------------------------------------------------------------
$ws7 = $wb.Worksheets.Item(7)

# Header row
$ws7.Cells.Item(1,2).Value  = "company"
$ws7.Cells.Item(1,3).Value  = "name"
$ws7.Cells.Item(1,4).Value  = "owner"
$ws7.Cells.Item(1,5).Value  = "property_category"
$ws7.Cells.Item(1,6).Value  = "category"
$ws7.Cells.Item(1,7).Value  = "date"
$ws7.Cells.Item(1,8).Value  = "legislator_name"
$ws7.Cells.Item(1,9).Value  = "legislator_id"
$ws7.Cells.Item(1,10).Value = "source_file"
$ws7.Cells.Item(1,11).Value = "index"

$insurance7 = @(
    @(85, "美國全民壽險公司", "20年期養老儲蓄險", "蔣乃辛"),
    @(86, "美國全民壽險公司", "20年期養老儲蓄險", "蔣乃辛"),
    @(87, "國泰人壽", "儲蓄型壽險", "楊際英"),
    @(88, "國泰人壽", "儲蓄型壽險", "楊際英"),
    @(89, "國泰人壽", "創世紀投資型壽險", "楊際英"),
    @(90, "國泰人壽", "創世紀投資型壽險", "楊際英"),
    @(91, "富邦人壽(原安泰ING)", "還本終身壽險", "楊際英"),
    @(92, "富邦人壽(原安泰ING)", "還本終身壽險", "楊際英")
)

$r = 2
foreach ($row in $insurance7) {
    $idx = $row[0]
    $ws7.Cells.Item($r,1).Value  = $idx
    $ws7.Cells.Item($r,2).Value  = $row[1]
    $ws7.Cells.Item($r,3).Value  = $row[2]
    $ws7.Cells.Item($r,4).Value  = $row[3]
    $ws7.Cells.Item($r,5).Value  = "insurance"
    $ws7.Cells.Item($r,6).Value  = "normal"
    # Leading apostrophe forces Excel to store this ISO-looking date as
    # plain text (shared string) instead of auto-converting it to a date
    # serial number, matching the source data's original text "date" field.
    $ws7.Cells.Item($r,7).Value  = "'2011-11-21"
    $ws7.Cells.Item($r,8).Value  = "蔣乃辛"
    $ws7.Cells.Item($r,9).Value  = 1722
    $ws7.Cells.Item($r,10).Value = "tmp12421"
    $ws7.Cells.Item($r,11).Value = $idx
    $r = $r + 1
}

# ---------------------------------------------------------------------
# Sheet 8: 債權 (claim)
# ---------------------------------------------------------------------
$ws8 = $wb.Worksheets.Item(8)

$ws8.Cells.Item(1,2).Value  = "species"
$ws8.Cells.Item(1,3).Value  = "owner"
$ws8.Cells.Item(1,4).Value  = "debtor"
$ws8.Cells.Item(1,5).Value  = "total"
$ws8.Cells.Item(1,6).Value  = "register_date"
$ws8.Cells.Item(1,7).Value  = "register_reason"
$ws8.Cells.Item(1,8).Value  = "property_category"
$ws8.Cells.Item(1,9).Value  = "category"
$ws8.Cells.Item(1,10).Value = "date"
$ws8.Cells.Item(1,11).Value = "legislator_name"
$ws8.Cells.Item(1,12).Value = "legislator_id"
$ws8.Cells.Item(1,13).Value = "source_file"
$ws8.Cells.Item(1,14).Value = "index"

$ws8.Cells.Item(2,1).Value  = 97
$ws8.Cells.Item(2,2).Value  = "暫付款"
$ws8.Cells.Item(2,3).Value  = "蔣乃辛"
$ws8.Cells.Item(2,4).Value  = "高美美臺北市羅斯福路"
$ws8.Cells.Item(2,5).Value  = 500000
$ws8.Cells.Item(2,6).Value  = 94
$ws8.Cells.Item(2,7).Value  = "借款（無法追回已被倒債）"
$ws8.Cells.Item(2,8).Value  = "claim"
$ws8.Cells.Item(2,9).Value  = "normal"
$ws8.Cells.Item(2,10).Value = "'2011-11-21"
$ws8.Cells.Item(2,11).Value = "蔣乃辛"
$ws8.Cells.Item(2,12).Value = 1722
$ws8.Cells.Item(2,13).Value = "tmp12421"
$ws8.Cells.Item(2,14).Value = 97

# ---------------------------------------------------------------------
# Sheet 9: 債務 (debt)
# ---------------------------------------------------------------------
$ws9 = $wb.Worksheets.Item(9)

$ws9.Cells.Item(1,2).Value  = "species"
$ws9.Cells.Item(1,3).Value  = "debtor"
$ws9.Cells.Item(1,4).Value  = "owner"
$ws9.Cells.Item(1,5).Value  = "total"
$ws9.Cells.Item(1,6).Value  = "register_date"
$ws9.Cells.Item(1,7).Value  = "register_reason"
$ws9.Cells.Item(1,8).Value  = "property_category"
$ws9.Cells.Item(1,9).Value  = "category"
$ws9.Cells.Item(1,10).Value = "date"
$ws9.Cells.Item(1,11).Value = "legislator_name"
$ws9.Cells.Item(1,12).Value = "legislator_id"
$ws9.Cells.Item(1,13).Value = "source_file"
$ws9.Cells.Item(1,14).Value = "index"

$ws9.Cells.Item(2,1).Value  = 102
$ws9.Cells.Item(2,2).Value  = "房貸"
$ws9.Cells.Item(2,3).Value  = "楊際英"
$ws9.Cells.Item(2,4).Value  = "永豐銀行三重分行新北市三重區忠孝路"
$ws9.Cells.Item(2,5).Value  = 4251536
$ws9.Cells.Item(2,6).Value  = "97年06月02日"
$ws9.Cells.Item(2,7).Value  = "金山南路2段房屋購屋貸款"
$ws9.Cells.Item(2,8).Value  = "debt"
$ws9.Cells.Item(2,9).Value  = "normal"
$ws9.Cells.Item(2,10).Value = "'2011-11-21"
$ws9.Cells.Item(2,11).Value = "蔣乃辛"
$ws9.Cells.Item(2,12).Value = 1722
$ws9.Cells.Item(2,13).Value = "tmp12421"
$ws9.Cells.Item(2,14).Value = 102
